# Update the "hindi" column (column E) annotations for rows 2-16
# after the hindi files rectification.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 29.46
    3  = 55.5
    4  = 17.4
    5  = 67.29000000000001
    6  = 34.54
    7  = 41.85
    8  = 74.34999999999999
    9  = 81.09
    10 = 71.67
    11 = 80.73
    12 = 84.44
    13 = 68.34999999999999
    14 = 52.08
    15 = 52.33
    16 = 85.14
}

foreach ($row in $updates.Keys) {
    $ws.Range("E$row").Value = $updates[$row]
}
